$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 388.53845
$ws.Range("I11").Value = 388.53845
$ws.Range("K11").Value = 388.53845
$ws.Range("M11").Value = -248.53845
$ws.Range("H64").Value = 4418
$ws.Range("I64").Value = 4420
$ws.Range("K64").Value = 4420
$ws.Range("M64").Value = -4172
$ws.Range("H67").Value = 4418
$ws.Range("I67").Value = 4420
$ws.Range("K67").Value = 4420
$ws.Range("M67").Value = -3562
$ws.Range("H70").Value = 3110
$ws.Range("I70").Value = 2831.6667
$ws.Range("J70").Value = 3666.6667
$ws.Range("K70").Value = 8495.000100000001
$ws.Range("L70").Value = 11000.0001
$ws.Range("M70").Value = -8225.000100000001
$ws.Range("N70").Value = -11540.0001
$ws.Range("H73").Value = 3110
$ws.Range("I73").Value = 2831.6667
$ws.Range("J73").Value = 3666.6667
$ws.Range("K73").Value = 8495.000100000001
$ws.Range("L73").Value = 11000.0001
$ws.Range("M73").Value = -7559.000100000001
$ws.Range("N73").Value = -12872.0001
$ws.Range("H116").Value = 5083.3335
$ws.Range("I116").Value = 4500
$ws.Range("J116").Value = 6250
$ws.Range("K116").Value = 4500
$ws.Range("L116").Value = 6250
$ws.Range("M116").Value = -1058
$ws.Range("N116").Value = -13134

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3254.8572
$ws.Range("I74").Value = 3254.8572
$ws.Range("K74").Value = 3254.8572
$ws.Range("M74").Value = -2380.8572
$ws.Range("H77").Value = 3254.8572
$ws.Range("I77").Value = 3254.8572
$ws.Range("K77").Value = 16274.286
$ws.Range("M77").Value = -11906.286

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2073.182
$ws.Range("I86").Value = 2170.5
$ws.Range("K86").Value = 2170.5
$ws.Range("M86").Value = -1047.5
$ws.Range("H88").Value = 20000
$ws.Range("J88").Value = 20000
$ws.Range("L88").Value = 20000
$ws.Range("N88").Value = -20812
$ws.Range("H89").Value = 2073.182
$ws.Range("I89").Value = 2170.5
$ws.Range("K89").Value = 10852.5
$ws.Range("M89").Value = -5236.5
$ws.Range("H91").Value = 20000
$ws.Range("J91").Value = 20000
$ws.Range("L91").Value = 20000
$ws.Range("N91").Value = -22808
$ws.Range("H94").Value = 2985.125
$ws.Range("I94").Value = 2701.3333
$ws.Range("J94").Value = 3836.5
$ws.Range("K94").Value = 2701.3333
$ws.Range("L94").Value = 3836.5
$ws.Range("M94").Value = -2250.3333
$ws.Range("N94").Value = -4738.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1976.3334
$ws.Range("I122").Value = 1720
$ws.Range("J122").Value = 3258
$ws.Range("K122").Value = 5160
$ws.Range("L122").Value = 9774
$ws.Range("M122").Value = -2710
$ws.Range("N122").Value = -14674
$ws.Range("H132").Value = 4238.6
$ws.Range("J132").Value = 4564
$ws.Range("L132").Value = 13692
$ws.Range("N132").Value = -18752

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2422.7144
$ws.Range("J129").Value = 2683.6
$ws.Range("L129").Value = 8050.799999999999
$ws.Range("N129").Value = -18050.8

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 21931.2
$ws.Range("I21").Value = 106
$ws.Range("K21").Value = 106
$ws.Range("M21").Value = 67
$ws.Range("H30").Value = 21931.2
$ws.Range("I30").Value = 106
$ws.Range("K30").Value = 106
$ws.Range("M30").Value = -1
$ws.Range("H80").Value = 2003.3334
$ws.Range("I80").Value = 2003.3334
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2003.3334
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -1005.3334
$ws.Range("N80").Value = $null
$ws.Range("H83").Value = 2003.3334
$ws.Range("I83").Value = 2003.3334
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 10016.667
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -5024.666999999999
$ws.Range("N83").Value = $null
$ws.Range("H113").Value = 950
$ws.Range("I113").Value = 950
$ws.Range("K113").Value = 950
$ws.Range("M113").Value = 1220
$ws.Range("H122").Value = 2724.3572
$ws.Range("I122").Value = 2434.2
$ws.Range("K122").Value = 7302.599999999999
$ws.Range("M122").Value = -4852.599999999999
$ws.Range("H126").Value = 15002.8
$ws.Range("I126").Value = 11000
$ws.Range("K126").Value = 33000
$ws.Range("M126").Value = -30530

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 7252.25
$ws.Range("I25").Value = 7
$ws.Range("J25").Value = 9667.333000000001
$ws.Range("K25").Value = 7
$ws.Range("L25").Value = 9667.333000000001
$ws.Range("M25").Value = 223
$ws.Range("N25").Value = -10127.333
$ws.Range("H82").Value = 5750
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 5750
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 5750
$ws.Range("M82").Value = $null
$ws.Range("N82").Value = -6472
$ws.Range("H85").Value = 5750
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 5750
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 5750
$ws.Range("M85").Value = $null
$ws.Range("N85").Value = -8246

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 975.2857
$ws.Range("I81").Value = 975.2857
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 1950.5714
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -889.5714
$ws.Range("N81").Value = $null
$ws.Range("H84").Value = 975.2857
$ws.Range("I84").Value = 975.2857
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 9752.857
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -4448.857
$ws.Range("H96").Value = 2397.6667
$ws.Range("J96").Value = 2477.2
$ws.Range("L96").Value = 2477.2
$ws.Range("N96").Value = -5223.2
$ws.Range("H97").Value = 7000
$ws.Range("J97").Value = 7000
$ws.Range("L97").Value = 7000
$ws.Range("N97").Value = -8982
$ws.Range("H107").Value = 1004.7895
$ws.Range("I107").Value = 786.0909
$ws.Range("J107").Value = 1305.5
$ws.Range("K107").Value = 2358.2727
$ws.Range("L107").Value = 3916.5
$ws.Range("M107").Value = -438.2727
$ws.Range("N107").Value = -7756.5
